$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "37.117.97"
$ws.Range("E2").Value = "  +0.05%  "

# Row 3
$ws.Range("D3").Value = "2.047.17"
$ws.Range("E3").Value = "  -0.45%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  -0.09%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "248.20"
$ws.Range("E5").Value = "  -0.56%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.664"
$ws.Range("E6").Value = "  -1.23%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "59.33"
$ws.Range("E7").Value = "  +0.32%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "1.00"
$ws.Range("E8").Value = "  -0.01%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.383"
$ws.Range("E9").Value = "  +1.35%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0784"
$ws.Range("E10").Value = "  -2.56%  "

# Row 11
$ws.Range("E11").Value = "  +0.77%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "15.82"
$ws.Range("E12").Value = "  +4.47%  "

# Row 13
$ws.Range("D13").Value = "2.345.24"
$ws.Range("E13").Value = "  -0.52%  "

# Row 14
$ws.Range("E14").Value = "  +2.45%  "

# Row 15
$ws.Range("E15").Value = "  +8.27%  "

# Row 16
$ws.Range("D16").Value = "2.061.52"
$ws.Range("E16").Value = "  +0.26%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "18.21"
$ws.Range("E17").Value = "  +25.66%  "

# Row 18
$ws.Range("D18").Value = "37.094.65"
$ws.Range("E18").Value = "  -0.04%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "74.66"
$ws.Range("E19").Value = "  -0.06%  "

# Row 20
$ws.Range("D20").Value = "0.0₃0895"
$ws.Range("E20").Value = "  -2.58%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.33"
$ws.Range("E21").Value = "  -0.29%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "236.69"
$ws.Range("E22").Value = "  -0.74%  "

# Row 23
$ws.Range("E23").Value = "  +0.02%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.45"
$ws.Range("E24").Value = "  +0.64%  "

# Row 25
$ws.Range("B25").Value = "Monero"
$ws.Range("C25").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "168.87"
$ws.Range("E25").Value = "  -1.71%  "

# Row 26
$ws.Range("B26").Value = "PancakeSwap"
$ws.Range("C26").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.16"
$ws.Range("E26").Value = "  +7.68%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.38"
$ws.Range("E27").Value = "  +2.36%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "19.97"
$ws.Range("E28").Value = "  -0.95%  "

# Row 29
$ws.Range("E29").Value = "  -0.06%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.13"
$ws.Range("E30").Value = "  +4.94%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.76"
$ws.Range("E31").Value = "  +3.14%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.0625"
$ws.Range("E32").Value = "  -1.21%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.52"
$ws.Range("E33").Value = "  +2.42%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0896"
$ws.Range("E34").Value = "  +1.71%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.00"
$ws.Range("E35").Value = "  +0.03%  "

# Row 36
$ws.Range("E36").Value = "  -2.58%  "

# Row 38
$ws.Range("B38").Value = "Cronos"
$ws.Range("C38").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.106"
$ws.Range("E38").Value = "  -1.65%  "

# Row 39
$ws.Range("B39").Value = "TrustWalletToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.33"
$ws.Range("E39").Value = "  -1.01%  "

# Row 40
$ws.Range("E40").Value = "  +13.52%  "

# Row 41
$ws.Range("E41").Value = "  +15.84%  "

# Row 42
$ws.Range("B42").Value = "VeChain"
$ws.Range("C42").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0221"
$ws.Range("E42").Value = "  -1.55%  "

# Row 43
$ws.Range("B43").Value = "InjectiveProtocol"
$ws.Range("C43").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "17.39"
$ws.Range("E43").Value = "  -4.92%  "

# Row 44
$ws.Range("E44").Value = "  -1.37%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "95.79"
$ws.Range("E45").Value = "  -1.17%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.46"
$ws.Range("E46").Value = "  -1.49%  "

# Row 47
$ws.Range("E47").Value = "  -0.07%  "

# Row 48
$ws.Range("D48").Value = "1.278.94"
$ws.Range("E48").Value = "  -1.86%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "6.76"
$ws.Range("E49").Value = "  -1.86%  "

# Row 50
$ws.Range("D50").Value = "2.226.58"
$ws.Range("E50").Value = "  -0.87%  "

# Row 51
$ws.Range("E51").Value = "  +0.94%  "
